$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 312
$ws.Range("I2").Value = 285.14285
$ws.Range("K2").Value = 285.14285
$ws.Range("M2").Value = -172.14285
$ws.Range("H17").Value = 528009.4399999999
$ws.Range("J17").Value = 528009.4399999999
$ws.Range("L17").Value = 1584028.32
$ws.Range("N17").Value = -1584364.32
$ws.Range("H28").Value = 424.0909
$ws.Range("I28").Value = 420.4762
$ws.Range("K28").Value = 420.4762
$ws.Range("M28").Value = 64.52379999999999
$ws.Range("H39").Value = 1822
$ws.Range("I39").Value = 2313.111
$ws.Range("J39").Value = 1190.5714
$ws.Range("K39").Value = 6939.333
$ws.Range("L39").Value = 3571.7142
$ws.Range("M39").Value = -6643.333
$ws.Range("N39").Value = -4163.7142
$ws.Range("H134").Value = 117702.57
$ws.Range("J134").Value = 107403.336
$ws.Range("L134").Value = 107403.336
$ws.Range("N134").Value = -117543.336
$ws.Range("H138").Value = 2922.7058
$ws.Range("J138").Value = 3494.2222
$ws.Range("L138").Value = 10482.6666
$ws.Range("N138").Value = -20762.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7524.391
$ws.Range("J45").Value = 3996.6924
$ws.Range("L45").Value = 3996.6924
$ws.Range("N45").Value = -4750.6924
$ws.Range("H61").Value = 1703.9344
$ws.Range("I61").Value = 1571.46
$ws.Range("J61").Value = 2306.0908
$ws.Range("K61").Value = 1571.46
$ws.Range("L61").Value = 2306.0908
$ws.Range("M61").Value = -1359.46
$ws.Range("N61").Value = -2730.0908
$ws.Range("H74").Value = 1800.1034
$ws.Range("I74").Value = 1485.4584
$ws.Range("K74").Value = 1485.4584
$ws.Range("M74").Value = -611.4584
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 1800.1034
$ws.Range("I77").Value = 1485.4584
$ws.Range("K77").Value = 7427.291999999999
$ws.Range("M77").Value = -3059.291999999999
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H110").Value = 2666.4
$ws.Range("I110").Value = 1346.3572
$ws.Range("K110").Value = 1346.3572
$ws.Range("M110").Value = 698.6428000000001
$ws.Range("H124").Value = 60740.5
$ws.Range("J124").Value = 60740.5
$ws.Range("L124").Value = 60740.5
$ws.Range("N124").Value = -70560.5
$ws.Range("H132").Value = 1846.174
$ws.Range("I132").Value = 1846.174
$ws.Range("K132").Value = 5538.522
$ws.Range("M132").Value = -3008.522
$ws.Range("H136").Value = 1703.9344
$ws.Range("I136").Value = 1571.46
$ws.Range("J136").Value = 2306.0908
$ws.Range("K136").Value = 4714.38
$ws.Range("L136").Value = 6918.2724
$ws.Range("M136").Value = -2164.38
$ws.Range("N136").Value = -12018.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6882.3237
$ws.Range("I20").Value = 7327.148
$ws.Range("K20").Value = 7327.148
$ws.Range("M20").Value = -7080.148
$ws.Range("H56").Value = 18110
$ws.Range("J56").Value = 18110
$ws.Range("L56").Value = 18110
$ws.Range("N56").Value = -19588
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H76").Value = 28438
$ws.Range("J76").Value = 27657
$ws.Range("L76").Value = 27657
$ws.Range("N76").Value = -28287
$ws.Range("H79").Value = 28438
$ws.Range("J79").Value = 27657
$ws.Range("L79").Value = 27657
$ws.Range("N79").Value = -29841
$ws.Range("H134").Value = 1334.1818
$ws.Range("I134").Value = 1282.6
$ws.Range("K134").Value = 3847.8
$ws.Range("M134").Value = -1312.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4774
$ws.Range("I31").Value = 3404.389
$ws.Range("K31").Value = 3404.389
$ws.Range("M31").Value = -3109.389
$ws.Range("H34").Value = 4774
$ws.Range("I34").Value = 3404.389
$ws.Range("K34").Value = 3404.389
$ws.Range("M34").Value = -3202.389
$ws.Range("H75").Value = 25260
$ws.Range("J75").Value = 25260
$ws.Range("L75").Value = 25260
$ws.Range("N75").Value = -27256
$ws.Range("H78").Value = 25260
$ws.Range("J78").Value = 25260
$ws.Range("L78").Value = 75780
$ws.Range("N78").Value = -85764
$ws.Range("H132").Value = 1593.9778
$ws.Range("I132").Value = 1547.8334
$ws.Range("K132").Value = 4643.5002
$ws.Range("M132").Value = -2113.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 10504840
$ws.Range("J55").Value = 17863098
$ws.Range("L55").Value = 53589294
$ws.Range("N55").Value = -53589648
$ws.Range("H80").Value = 4142.1665
$ws.Range("J80").Value = 4275.7334
$ws.Range("L80").Value = 12827.2002
$ws.Range("N80").Value = -14699.2002
$ws.Range("H83").Value = 4142.1665
$ws.Range("J83").Value = 4275.7334
$ws.Range("L83").Value = 38481.6006
$ws.Range("N83").Value = -47841.6006
$ws.Range("H103").Value = 766.3333
$ws.Range("I103").Value = 466.33334
$ws.Range("J103").Value = 1066.3334
$ws.Range("K103").Value = 1399.00002
$ws.Range("L103").Value = 3199.0002
$ws.Range("M103").Value = -520.0000199999999
$ws.Range("N103").Value = -4957.0002
$ws.Range("H122").Value = 725.9231
$ws.Range("I122").Value = 723.8
$ws.Range("K122").Value = 6514.2
$ws.Range("M122").Value = -4064.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2562
$ws.Range("J102").Value = 2199.25
$ws.Range("L102").Value = 2199.25
$ws.Range("N102").Value = -5443.25
$ws.Range("H123").Value = 67081.60000000001
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 67081.60000000001
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 67081.60000000001
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -71981.60000000001
$ws.Range("H126").Value = 2827.6
$ws.Range("I126").Value = 3120.3333
$ws.Range("J126").Value = 2388.5
$ws.Range("K126").Value = 9360.999899999999
$ws.Range("L126").Value = 7165.5
$ws.Range("M126").Value = -6890.999899999999
$ws.Range("N126").Value = -12105.5
$ws.Range("H132").Value = 6373.4
$ws.Range("I132").Value = 4916.0713
$ws.Range("K132").Value = 14748.2139
$ws.Range("M132").Value = -12218.2139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15752.964
$ws.Range("I40").Value = 11220.929
$ws.Range("J40").Value = 20285
$ws.Range("K40").Value = 11220.929
$ws.Range("L40").Value = 20285
$ws.Range("M40").Value = -11084.929
$ws.Range("N40").Value = -20557
$ws.Range("H122").Value = 83668.8
$ws.Range("I122").Value = 108160.734
$ws.Range("K122").Value = 324482.202
$ws.Range("M122").Value = -322032.202

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1217799.8
$ws.Range("I14").Value = 3001999.5
$ws.Range("K14").Value = 3001999.5
$ws.Range("M14").Value = -3001831.5
$ws.Range("H113").Value = 365.66666
$ws.Range("I113").Value = 338.625
$ws.Range("K113").Value = 1015.875
$ws.Range("M113").Value = 1154.125
$ws.Range("H132").Value = 1061101.6
$ws.Range("I132").Value = 2383.9688
$ws.Range("J132").Value = 3667176
$ws.Range("K132").Value = 7151.9064
$ws.Range("L132").Value = 11001528
$ws.Range("M132").Value = -4621.9064
$ws.Range("N132").Value = -11006588
$ws.Range("H136").Value = 970.7
$ws.Range("I136").Value = 970.7
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2912.1
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -362.1000000000004
$ws.Range("N136").ClearContents()
